$wb = $excel.ActiveWorkbook

# --- Sheet "ANLT" (first sheet) ---
$ws1 = $wb.Worksheets.Item("ANLT")

# --- Sheet "ANHDT" (second sheet) ---
$ws2 = $wb.Worksheets.Item("ANHDT")

# Shared-string allocation order must match the authored order:
# 27: ws1!A11, 28: ws2!A4, 29: ws2!A5, 30: ws1!A12
$ws1.Range("A11").Value = "Update cơ chế phần đa ngôn ngữ cho sub và update ngôn ngữ trên home"
$ws1.Range("M11").Value = 4
$ws1.Range("N11").Value = 0

$ws2.Range("A4").Value = "Details chưa có rating & review do đã thống nhất chưa làm. Hiện tại đã có thể bấm vào add to cart và add to wishlist"
$ws2.Range("A4").WrapText = $true
$ws2.Rows.Item(4).RowHeight = 45
$ws2.Range("N4").Value = 4

$ws2.Range("A5").Value = "Metting"
$ws2.Range("O5").Value = 4

$ws1.Range("A12").Value = "Metting + Support(Time và task dự án)"
$ws1.Range("O12").Value = 4

# --- Selections / active sheet ---
$ws2.Range("O15").Select()

$ws1.Activate()
$ws1.Range("B19").Select()
